$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: repeat of the header row (row 1) -----------------------------
# Copying preserves both the values and the cell formatting/styles (and
# reuses existing shared-string entries instead of creating duplicates).
$ws.Range("A1:E1").Copy($ws.Range("A11:E11"))
$ws.Range("G1:J1").Copy($ws.Range("G11:J11"))

# --- Row 12: "Lloyd_slow_seq" timings --------------------------------------
$ws.Range("A12").Value = "Lloyd_slow_seq"
$ws.Range("B12").Value = 6.41404
$ws.Range("C12").Value = 18.559297
$ws.Range("D12").Value = 100.31866
$ws.Range("E12").Value = 307.145723

# --- Row 13: "lloyd_sl_par" timings ----------------------------------------
$ws.Range("A13").Value = "lloyd_sl_par"
$ws.Range("B13").Value = 2.65987
$ws.Range("C13").Value = 7.356104
$ws.Range("D13").Value = 42.23705
$ws.Range("E13").Value = 131.894881

# --- Row 14: speedup = row12 / row13 ---------------------------------------
$ws.Range("B14").Formula = "=B12/B13"
$ws.Range("C14").Formula = "=C12/C13"
$ws.Range("D14").Formula = "=D12/D13"
$ws.Range("E14").Formula = "=E12/E13"

# --- Update the sheet view: scroll down and change the active selection ---
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C31").Select() | Out-Null
